# Auto-generated script applying scheduled market-data refresh to Carbuncle_Profits workbook.
# For each sheet, update the currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# with freshly polled values. Cells that had no prior value (blank) are populated; the one
# cell whose value is now redundant (superseded by its neighbour) is cleared to blank.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1448.625
$ws.Range("I40").Value = 1616
$ws.Range("J40").Value = 1080.4
$ws.Range("K40").Value = 1616
$ws.Range("L40").Value = 1080.4
$ws.Range("M40").Value = -1441
$ws.Range("N40").Value = -1430.4
$ws.Range("H98").Value = 1425.1428
$ws.Range("I98").Value = 1596
$ws.Range("J98").Value = 998
$ws.Range("K98").Value = 1596
$ws.Range("L98").Value = 998
$ws.Range("M98").Value = -98
$ws.Range("N98").Value = -3994
$ws.Range("H122").Value = 1425.1428
$ws.Range("I122").Value = 1596
$ws.Range("J122").Value = 998
$ws.Range("K122").Value = 4788
$ws.Range("L122").Value = 2994
$ws.Range("M122").Value = -2338
$ws.Range("N122").Value = -7894
$ws.Range("H132").Value = 1875.9445
$ws.Range("I132").Value = 1996.0312
$ws.Range("K132").Value = 5988.0936
$ws.Range("M132").Value = -3458.0936

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3397.2273
$ws.Range("I32").Value = 2174.6753
$ws.Range("K32").Value = 2174.6753
$ws.Range("M32").Value = -1887.6753
$ws.Range("H122").Value = 83336090
$ws.Range("I122").Value = 111113460
$ws.Range("J122").Value = 3999
$ws.Range("K122").Value = 333340380
$ws.Range("L122").Value = 11997
$ws.Range("M122").Value = -333337930
$ws.Range("N122").Value = -16897

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 7300
$ws.Range("I75").Value = 4050
$ws.Range("K75").Value = 4050
$ws.Range("M75").Value = -3114
$ws.Range("H78").Value = 7300
$ws.Range("I78").Value = 4050
$ws.Range("K78").Value = 12150
$ws.Range("M78").Value = -7470
$ws.Range("H99").Value = 1323.7073
$ws.Range("I99").Value = 797.2632
$ws.Range("K99").Value = 797.2632
$ws.Range("M99").Value = 700.7368
$ws.Range("H134").Value = 1934.4348
$ws.Range("I134").Value = 942.375
$ws.Range("K134").Value = 2827.125
$ws.Range("M134").Value = -292.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5352570.5
$ws.Range("I99").Value = 8003231
$ws.Range("J99").Value = 51249.5
$ws.Range("K99").Value = 8003231
$ws.Range("L99").Value = 51249.5
$ws.Range("M99").Value = -8001733
$ws.Range("N99").Value = -54245.5
$ws.Range("H125").Value = 26889.875
$ws.Range("J125").Value = 26889.875
$ws.Range("L125").Value = 26889.875
$ws.Range("N125").Value = -31809.875
$ws.Range("H126").Value = 5352570.5
$ws.Range("I126").Value = 8003231
$ws.Range("J126").Value = 51249.5
$ws.Range("K126").Value = 24009693
$ws.Range("L126").Value = 153748.5
$ws.Range("M126").Value = -24007223
$ws.Range("N126").Value = -158688.5
$ws.Range("H132").Value = 31490.5
$ws.Range("I132").Value = 37457.68
$ws.Range("J132").Value = 3643.6667
$ws.Range("K132").Value = 112373.04
$ws.Range("L132").Value = 10931.0001
$ws.Range("M132").Value = -109843.04
$ws.Range("N132").Value = -15991.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1353.7693
$ws.Range("I4").Value = 359.8
$ws.Range("K4").Value = 1079.4
$ws.Range("M4").Value = -967.4000000000001
$ws.Range("H5").Value = 283804.88
$ws.Range("I5").Value = 370.04166
$ws.Range("J5").Value = 462816.38
$ws.Range("K5").Value = 1110.12498
$ws.Range("L5").Value = 1388449.14
$ws.Range("M5").Value = -998.1249800000001
$ws.Range("N5").Value = -1388673.14
$ws.Range("H35").Value = 837.5
$ws.Range("I35").Value = 2300
$ws.Range("J35").Value = 350
$ws.Range("K35").Value = 6900
$ws.Range("L35").Value = 1050
$ws.Range("M35").Value = -6612
$ws.Range("N35").Value = -1626
$ws.Range("H49").Value = 15000
$ws.Range("J49").Value = 15000
$ws.Range("L49").Value = 45000
$ws.Range("N49").Value = -45312
$ws.Range("H57").Value = 8000
$ws.Range("J57").Value = 8000
$ws.Range("L57").Value = 24000
$ws.Range("N57").Value = -25118
$ws.Range("H68").Value = 564078.9
$ws.Range("I68").Value = 1735932.2
$ws.Range("J68").Value = 1589.22
$ws.Range("K68").Value = 5207796.6
$ws.Range("L68").Value = 4767.66
$ws.Range("M68").Value = -5206985.6
$ws.Range("N68").Value = -6389.66
$ws.Range("H71").Value = 564078.9
$ws.Range("I71").Value = 1735932.2
$ws.Range("J71").Value = 1589.22
$ws.Range("K71").Value = 15623389.8
$ws.Range("L71").Value = 14302.98
$ws.Range("M71").Value = -15619333.8
$ws.Range("N71").Value = -22414.98
$ws.Range("H74").Value = 2999
$ws.Range("J74").Value = 2999
$ws.Range("L74").Value = 8997
$ws.Range("N74").Value = -11119
$ws.Range("H77").Value = 2999
$ws.Range("J77").Value = 2999
$ws.Range("L77").Value = 26991
$ws.Range("N77").Value = -37599
$ws.Range("H94").Value = 5400
$ws.Range("I94").Value = 3000
$ws.Range("J94").Value = 6200
$ws.Range("K94").Value = 9000
$ws.Range("L94").Value = 18600
$ws.Range("M94").Value = -8324
$ws.Range("N94").Value = -19952
$ws.Range("H99").Value = 2100
$ws.Range("I99").Value = 500
$ws.Range("J99").Value = 2900
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 8700
$ws.Range("M99").Value = 746
$ws.Range("N99").Value = -13192
$ws.Range("H101").Value = 14996.333
$ws.Range("J101").Value = 14996.333
$ws.Range("L101").Value = 44988.999
$ws.Range("N101").Value = -49856.999
$ws.Range("H122").Value = 599.7241
$ws.Range("I122").Value = 391.76
$ws.Range("J122").Value = 1899.5
$ws.Range("K122").Value = 3525.84
$ws.Range("L122").Value = 17095.5
$ws.Range("M122").Value = -1075.84
$ws.Range("N122").Value = -21995.5
$ws.Range("H131").Value = 5509.68
$ws.Range("I131").Value = 1100
$ws.Range("J131").Value = 5693.4165
$ws.Range("K131").Value = 3300
$ws.Range("L131").Value = 17080.2495
$ws.Range("M131").Value = 1740
$ws.Range("N131").Value = -27160.2495
$ws.Range("H132").Value = 1234.6562
$ws.Range("I132").Value = 1390.9445
$ws.Range("J132").Value = 1033.7142
$ws.Range("K132").Value = 12518.5005
$ws.Range("L132").Value = 9303.427799999999
$ws.Range("M132").Value = -9988.5005
$ws.Range("N132").Value = -14363.4278
$ws.Range("H133").Value = 1891.6666
$ws.Range("H134").Value = 38436.258
$ws.Range("J134").Value = 2027.2727
$ws.Range("L134").Value = 6081.8181
$ws.Range("N134").Value = -16221.8181
$ws.Range("H135").Value = 283804.88
$ws.Range("I135").Value = 370.04166
$ws.Range("J135").Value = 462816.38
$ws.Range("K135").Value = 3330.37494
$ws.Range("L135").Value = 4165347.42
$ws.Range("M135").Value = -795.3749399999997
$ws.Range("N135").Value = -4170417.42
$ws.Range("H137").Value = 2012.2593
$ws.Range("I137").Value = 1956.125
$ws.Range("J137").Value = 2093.9092
$ws.Range("K137").Value = 5868.375
$ws.Range("L137").Value = 6281.7276
$ws.Range("M137").Value = -768.375
$ws.Range("N137").Value = -16481.7276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14737.889
$ws.Range("I70").Value = 20140.5
$ws.Range("K70").Value = 20140.5
$ws.Range("M70").Value = -19870.5
$ws.Range("H73").Value = 14737.889
$ws.Range("I73").Value = 20140.5
$ws.Range("K73").Value = 20140.5
$ws.Range("M73").Value = -19204.5
$ws.Range("H80").Value = 2166.5833
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 2399.8
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 2399.8
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -4395.8
$ws.Range("H83").Value = 2166.5833
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 2399.8
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 11999
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -21983

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 62423.35
$ws.Range("I7").Value = 80553.62
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 80553.62
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -80441.62
$ws.Range("N7").Value = -3724
$ws.Range("H93").Value = 14969
$ws.Range("I93").Value = 20608.6
$ws.Range("J93").Value = 870
$ws.Range("K93").Value = 20608.6
$ws.Range("L93").Value = 870
$ws.Range("M93").Value = -19360.6
$ws.Range("N93").Value = -3366
$ws.Range("H122").Value = 111111110
$ws.Range("I122").Value = 111111110
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 333333330
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -333330880
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 62423.35
$ws.Range("I126").Value = 80553.62
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 241660.86
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -239190.86
$ws.Range("N126").Value = -15440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 36658.668
$ws.Range("I75").Value = 20000
$ws.Range("K75").Value = 20000
$ws.Range("M75").Value = -19064
$ws.Range("H78").Value = 36658.668
$ws.Range("I78").Value = 20000
$ws.Range("K78").Value = 60000
$ws.Range("M78").Value = -55320
$ws.Range("H81").Value = 1373.5333
$ws.Range("I81").Value = 1787.625
$ws.Range("K81").Value = 3575.25
$ws.Range("M81").Value = -2514.25
$ws.Range("H84").Value = 1373.5333
$ws.Range("I84").Value = 1787.625
$ws.Range("K84").Value = 17876.25
$ws.Range("M84").Value = -12572.25
$ws.Range("H132").Value = 2243.725
$ws.Range("I132").Value = 2309.5715
$ws.Range("J132").Value = 2208.2693
$ws.Range("K132").Value = 6928.7145
$ws.Range("L132").Value = 6624.8079
$ws.Range("M132").Value = -4398.7145
$ws.Range("N132").Value = -11684.8079
